# Update gh-pages to output generated at 456a3b4
#
# Applies the edits described by the diff to both the "展览" sheet
# (Worksheets index 1) and the "全部类型" sheet (Worksheets index 4) -
# they list (mostly) the same events, merged with other categories on
# the "全部类型" sheet, so matching rows live at different row numbers
# on each sheet. Row numbers below are therefore spelled out per sheet
# rather than derived from a shared offset.

$wb = $excel.ActiveWorkbook

function Set-TextValue($ws, [string]$addr, [string]$text) {
    # Force a plain text cell even when the text looks like a date
    # (e.g. "2024-10-04") so Excel doesn't silently convert it to a
    # date serial number. Clearing formats afterwards drops the
    # temporary "@" number format so the cell matches a freshly
    # authored text cell (no residual style index).
    $ws.Range($addr).NumberFormat = "@"
    $ws.Range($addr).Value = $text
    $ws.Range($addr).ClearFormats()
}

function Set-NewRow($ws, [int]$styleSrcRow, [int]$rowNum, [int]$a, [string]$b, [string]$c, [string]$d, [string]$e, $f, $g, [string]$h, [string]$i) {
    # Clone the "serial number" column's look (bold, centered, bordered)
    # from another data row instead of re-building it property by
    # property - that keeps the exact same style index as the existing
    # rows instead of registering a near-duplicate style.
    $ws.Range("A$styleSrcRow").Copy()
    $ws.Range("A$rowNum").PasteSpecial(-4122)
    $ws.Range("A$rowNum").Value = $a

    Set-TextValue $ws "B$rowNum" $b
    $ws.Range("C$rowNum").Value = $c
    $ws.Range("D$rowNum").Value = $d
    $ws.Range("E$rowNum").Value = $e
    $ws.Range("F$rowNum").Value = $f
    $ws.Range("G$rowNum").Value = $g
    $ws.Range("H$rowNum").Value = $h
    $ws.Range("I$rowNum").Value = $i
}

# ============================= Sheet "展览" =============================
$ws1 = $wb.Worksheets.Item(1)

$ws1.Range("F2").Value = 230
$ws1.Range("F3").Value = 1076

$ws1.Range("C4").Value = "张家港·META萌圆饿了（取消）"
$ws1.Range("G4").Value = "不可售"

$ws1.Range("F8").Value = 50
$ws1.Range("F9").Value = 6683
$ws1.Range("F10").Value = 141

$ws1.Range("G12").Value = "不可售"

$ws1.Range("F16").Value = 15988
$ws1.Range("F18").Value = 29

$ws1.Range("F22").Value = 11244
$ws1.Range("F23").Value = 826
$ws1.Range("F24").Value = 4419
$ws1.Range("F25").Value = 288

# Insert a brand-new row at 28 ("常熟·cc动漫游戏嘉年华"), pushing the
# existing rows 28-29 down to 29-30.
$ws1.Rows.Item(28).Insert()
Set-NewRow $ws1 27 28 27 "2024-10-04" "常熟·cc动漫游戏嘉年华" "开元大道1号 常熟国际博览中心" "2024.10.04 09:00-10.05 17:00" 0 60 "https://show.bilibili.com/platform/detail.html?id=90292" "//i2.hdslb.com/bfs/openplatform/202407/yCNXedrA1722404050722.jpeg"

# Fix up the A (serial #) column for the two rows that shifted down.
$ws1.Range("A29").Value = 28
$ws1.Range("A30").Value = 29

# Append a brand-new row 31 ("苏州·星部落动漫嘉年华").
Set-NewRow $ws1 30 31 30 "2024-12-27" "苏州·星部落动漫嘉年华" "花桥经济开发区绿地大道1598号 花桥国际博览中心" "2024.12.27 09:00-12.28 16:00" 5213 68 "https://show.bilibili.com/platform/detail.html?id=84858" "//i0.hdslb.com/bfs/openplatform/202404/UI5EFZTT1713685680462.jpeg"

# ============================ Sheet "全部类型" ============================
$ws4 = $wb.Worksheets.Item(4)

$ws4.Range("F2").Value = 230
$ws4.Range("F3").Value = 1076

$ws4.Range("C4").Value = "张家港·META萌圆饿了（取消）"
$ws4.Range("G4").Value = "不可售"

$ws4.Range("F9").Value = 50
$ws4.Range("F10").Value = 6683
$ws4.Range("F11").Value = 141

$ws4.Range("G13").Value = "不可售"

$ws4.Range("F18").Value = 15988
$ws4.Range("F20").Value = 29

$ws4.Range("F25").Value = 11245
$ws4.Range("F26").Value = 826
$ws4.Range("F27").Value = 4419
$ws4.Range("F28").Value = 288

# Insert a brand-new row at 31 ("常熟·cc动漫游戏嘉年华"), pushing the
# existing rows 31-32 down to 32-33.
$ws4.Rows.Item(31).Insert()
Set-NewRow $ws4 30 31 30 "2024-10-04" "常熟·cc动漫游戏嘉年华" "开元大道1号 常熟国际博览中心" "2024.10.04 09:00-10.05 17:00" 0 60 "https://show.bilibili.com/platform/detail.html?id=90292" "//i2.hdslb.com/bfs/openplatform/202407/yCNXedrA1722404050722.jpeg"

# Fix up the A (serial #) column for the two rows that shifted down.
$ws4.Range("A32").Value = 31
$ws4.Range("A33").Value = 32

# Append a brand-new row 34 ("苏州·星部落动漫嘉年华").
Set-NewRow $ws4 33 34 33 "2024-12-27" "苏州·星部落动漫嘉年华" "花桥经济开发区绿地大道1598号 花桥国际博览中心" "2024.12.27 09:00-12.28 16:00" 5213 68 "https://show.bilibili.com/platform/detail.html?id=84858" "//i0.hdslb.com/bfs/openplatform/202404/UI5EFZTT1713685680462.jpeg"
